$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.922.46'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').Value = '1.744.47'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.13%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '247.12'
$ws.Range('E5').Value = '  +3.47%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.9998'
$ws.Range('E6').Value = '  +0.15%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.5042'
$ws.Range('E7').Value = '  -4.89%  '
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '40.77'
$ws.Range('E8').Value = '  +1.36%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.2736'
$ws.Range('E9').Value = '  -3.20%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.06187'
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '1.751.53'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.07257'
$ws.Range('E12').Value = '  +0.63%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.6534'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '15.15'
$ws.Range('E14').Value = '  -3.11%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '4.635'
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '77.53'
$ws.Range('E16').Value = '  -1.85%  '
$ws.Range('B17').Value = 'Dai'
$ws.Range('C17').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '1.000'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.9995'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '25.945.16'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '11.84'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '0.000006820'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '1.974.95'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '4.346'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '8.674'
$ws.Range('E24').Value = '  -1.44%  '
$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '5.408'
$ws.Range('E25').Value = '  +2.57%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '136.90'
$ws.Range('E26').Value = '  -1.92%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '1.499'
$ws.Range('E27').Value = '  -2.16%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '15.20'
$ws.Range('E28').Value = '  -1.47%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.772'
$ws.Range('E29').Value = '  -3.18%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '105.72'
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '3.914'
$ws.Range('E31').Value = '  +1.92%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.08233'
$ws.Range('E32').Value = '  -1.20%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '3.623'
$ws.Range('E33').Value = '  -1.47%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.04675'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '2.655'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.9954'
$ws.Range('E36').Value = '  -3.20%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.6187'
$ws.Range('E37').Value = '  -3.45%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '2.738'
$ws.Range('E38').Value = '  +1.09%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.01607'
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '1.916'
$ws.Range('E40').Value = '  -4.44%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.9996'
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '99.90'
$ws.Range('E42').Value = '  -2.57%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.3888'
$ws.Range('E43').Value = '  -2.30%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.7569'
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '5.000'
$ws.Range('E45').Value = '  -1.19%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.1146'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '6.310'
$ws.Range('E47').Value = '  -2.12%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '55.51'
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.05249'
$ws.Range('E49').Value = '  -1.91%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '30.60'
$ws.Range('E50').Value = '  -1.79%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '7.520'
$ws.Range('E51').Value = '  -1.36%  '

Write-Host "Applied all changes"
